$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7: Experimental -> set Value to literal text "false" (not boolean FALSE).
# A direct $cell.Value = "false" assignment gets auto-coerced by Excel into a
# Boolean cell; going through a text formula + paste-as-values keeps it a
# genuine string cell like the source FHIR IG generator produced.
$cellExperimental = $ws.Range("B7")
$cellExperimental.Formula = "=""false"""
$cellExperimental.Copy()
$cellExperimental.PasteSpecial(-4163)  # xlPasteValues

# Row 8: Date -> update the ISO timestamp value.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Row 17: Description -> set the description value (was empty).
$ws.Range("B17").Value = "Codes for accumulated recovery debt levels"

$excel.CutCopyMode = 0
